$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "eclipse-error-the-import-xxx-cannot-be-resolved"
$ws.Range("B4").Value = "Go to project in the task view and click on Clean.."

$ws.Range("B4").Select()
